# Auto-generated edit script: updates live market-price-derived columns
# (currentAveragePrice* / LevePrice* / LeveProfit*) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 117.71429
$ws.Cells.Item(5, 9).Value = 78.916664
$ws.Cells.Item(5, 10).Value = 350.5
$ws.Cells.Item(5, 11).Value = 78.916664
$ws.Cells.Item(5, 12).Value = 350.5
$ws.Cells.Item(5, 13).Value = 36.083336
$ws.Cells.Item(5, 14).Value = -580.5

$ws.Cells.Item(51, 8).Value = 5740.5
$ws.Cells.Item(51, 9).Value = 4071.4285
$ws.Cells.Item(51, 10).Value = 7687.75
$ws.Cells.Item(51, 11).Value = 4071.4285
$ws.Cells.Item(51, 12).Value = 7687.75
$ws.Cells.Item(51, 13).Value = -3587.4285
$ws.Cells.Item(51, 14).Value = -8655.75

$ws.Cells.Item(96, 8).Value = 1185.1875
$ws.Cells.Item(96, 9).Value = 372.42856
$ws.Cells.Item(96, 10).Value = 1817.3334
$ws.Cells.Item(96, 11).Value = 1117.28568
$ws.Cells.Item(96, 12).Value = 5452.0002
$ws.Cells.Item(96, 13).Value = 255.71432
$ws.Cells.Item(96, 14).Value = -8198.0002

$ws.Cells.Item(98, 8).Value = 3133.7144
$ws.Cells.Item(98, 9).Value = 2598.4
$ws.Cells.Item(98, 10).Value = 4472
$ws.Cells.Item(98, 11).Value = 2598.4
$ws.Cells.Item(98, 12).Value = 4472
$ws.Cells.Item(98, 13).Value = -1100.4
$ws.Cells.Item(98, 14).Value = -7468

$ws.Cells.Item(113, 8).Value = 2929.2188
$ws.Cells.Item(113, 9).Value = 2678.2856
$ws.Cells.Item(113, 10).Value = 3124.389
$ws.Cells.Item(113, 11).Value = 2678.2856
$ws.Cells.Item(113, 12).Value = 3124.389
$ws.Cells.Item(113, 13).Value = 575.7143999999998
$ws.Cells.Item(113, 14).Value = -9632.388999999999

$ws.Cells.Item(122, 8).Value = 3133.7144
$ws.Cells.Item(122, 9).Value = 2598.4
$ws.Cells.Item(122, 10).Value = 4472
$ws.Cells.Item(122, 11).Value = 7795.200000000001
$ws.Cells.Item(122, 12).Value = 13416
$ws.Cells.Item(122, 13).Value = -5345.200000000001
$ws.Cells.Item(122, 14).Value = -18316

$ws.Cells.Item(132, 8).Value = 98367.63
$ws.Cells.Item(132, 9).Value = 211083.42
$ws.Cells.Item(132, 10).Value = 13336.421
$ws.Cells.Item(132, 11).Value = 633250.26
$ws.Cells.Item(132, 12).Value = 40009.263
$ws.Cells.Item(132, 13).Value = -630720.26
$ws.Cells.Item(132, 14).Value = -45069.263

$ws.Cells.Item(137, 8).Value = 4517.6
$ws.Cells.Item(137, 9).Value = 1900
$ws.Cells.Item(137, 10).Value = 4808.4443
$ws.Cells.Item(137, 11).Value = 5700
$ws.Cells.Item(137, 12).Value = 14425.3329
$ws.Cells.Item(137, 13).Value = -3150
$ws.Cells.Item(137, 14).Value = -19525.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4840.391
$ws.Cells.Item(32, 9).Value = 2791.0286
$ws.Cells.Item(32, 10).Value = 11361.091
$ws.Cells.Item(32, 11).Value = 2791.0286
$ws.Cells.Item(32, 12).Value = 11361.091
$ws.Cells.Item(32, 13).Value = -2504.0286
$ws.Cells.Item(32, 14).Value = -11935.091

$ws.Cells.Item(43, 8).Value = 20000
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 20000
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 20000
$ws.Cells.Item(43, 14).Value = -20626
$ws.Cells.Item(43, 13).ClearContents()

$ws.Cells.Item(45, 8).Value = 2376.2
$ws.Cells.Item(45, 9).Value = 2376.2
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 2376.2
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = -1999.2

$ws.Cells.Item(61, 8).Value = 12281.806
$ws.Cells.Item(61, 9).Value = 9495.478999999999
$ws.Cells.Item(61, 10).Value = 17211.46
$ws.Cells.Item(61, 11).Value = 9495.478999999999
$ws.Cells.Item(61, 12).Value = 17211.46
$ws.Cells.Item(61, 13).Value = -9283.478999999999
$ws.Cells.Item(61, 14).Value = -17635.46

$ws.Cells.Item(74, 8).Value = 1253.2916
$ws.Cells.Item(74, 9).Value = 750.1818
$ws.Cells.Item(74, 10).Value = 1679
$ws.Cells.Item(74, 11).Value = 750.1818
$ws.Cells.Item(74, 12).Value = 1679
$ws.Cells.Item(74, 13).Value = 123.8182
$ws.Cells.Item(74, 14).Value = -3427

$ws.Cells.Item(77, 8).Value = 1253.2916
$ws.Cells.Item(77, 9).Value = 750.1818
$ws.Cells.Item(77, 10).Value = 1679
$ws.Cells.Item(77, 11).Value = 3750.909
$ws.Cells.Item(77, 12).Value = 8395
$ws.Cells.Item(77, 13).Value = 617.0910000000003
$ws.Cells.Item(77, 14).Value = -17131

$ws.Cells.Item(122, 8).Value = 4693.5713
$ws.Cells.Item(122, 9).Value = 2588.2144
$ws.Cells.Item(122, 10).Value = 8904.286
$ws.Cells.Item(122, 11).Value = 7764.6432
$ws.Cells.Item(122, 12).Value = 26712.858
$ws.Cells.Item(122, 13).Value = -5314.6432
$ws.Cells.Item(122, 14).Value = -31612.858

$ws.Cells.Item(132, 8).Value = 3547.6
$ws.Cells.Item(132, 9).Value = 1673.7941
$ws.Cells.Item(132, 10).Value = 14165.833
$ws.Cells.Item(132, 11).Value = 5021.3823
$ws.Cells.Item(132, 12).Value = 42497.499
$ws.Cells.Item(132, 13).Value = -2491.3823
$ws.Cells.Item(132, 14).Value = -47557.499

$ws.Cells.Item(136, 8).Value = 12281.806
$ws.Cells.Item(136, 9).Value = 9495.478999999999
$ws.Cells.Item(136, 10).Value = 17211.46
$ws.Cells.Item(136, 11).Value = 28486.437
$ws.Cells.Item(136, 12).Value = 51634.38
$ws.Cells.Item(136, 13).Value = -25936.437
$ws.Cells.Item(136, 14).Value = -56734.38

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4146.143
$ws.Cells.Item(20, 9).Value = 3337.1667
$ws.Cells.Item(20, 10).Value = 9000
$ws.Cells.Item(20, 11).Value = 3337.1667
$ws.Cells.Item(20, 12).Value = 9000
$ws.Cells.Item(20, 13).Value = -3090.1667
$ws.Cells.Item(20, 14).Value = -9494

$ws.Cells.Item(86, 8).Value = 4289.6
$ws.Cells.Item(86, 9).Value = 3964.6667
$ws.Cells.Item(86, 10).Value = 4777
$ws.Cells.Item(86, 11).Value = 3964.6667
$ws.Cells.Item(86, 12).Value = 4777
$ws.Cells.Item(86, 13).Value = -2841.6667
$ws.Cells.Item(86, 14).Value = -7023

$ws.Cells.Item(89, 8).Value = 4289.6
$ws.Cells.Item(89, 9).Value = 3964.6667
$ws.Cells.Item(89, 10).Value = 4777
$ws.Cells.Item(89, 11).Value = 19823.3335
$ws.Cells.Item(89, 12).Value = 23885
$ws.Cells.Item(89, 13).Value = -14207.3335
$ws.Cells.Item(89, 14).Value = -35117

$ws.Cells.Item(134, 8).Value = 3241.75
$ws.Cells.Item(134, 9).Value = 2767.2778
$ws.Cells.Item(134, 10).Value = 4665.1665
$ws.Cells.Item(134, 11).Value = 8301.8334
$ws.Cells.Item(134, 12).Value = 13995.4995
$ws.Cells.Item(134, 13).Value = -5766.8334
$ws.Cells.Item(134, 14).Value = -19065.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 4577.75
$ws.Cells.Item(122, 9).Value = 1739
$ws.Cells.Item(122, 10).Value = 7416.5
$ws.Cells.Item(122, 11).Value = 5217
$ws.Cells.Item(122, 12).Value = 22249.5
$ws.Cells.Item(122, 13).Value = -2767
$ws.Cells.Item(122, 14).Value = -27149.5

$ws.Cells.Item(132, 8).Value = 13187.5
$ws.Cells.Item(132, 9).Value = 4861.1113
$ws.Cells.Item(132, 10).Value = 38166.668
$ws.Cells.Item(132, 11).Value = 14583.3339
$ws.Cells.Item(132, 12).Value = 114500.004
$ws.Cells.Item(132, 13).Value = -12053.3339
$ws.Cells.Item(132, 14).Value = -119560.004

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 64129092
$ws.Cells.Item(4, 9).Value = 1615738.2
$ws.Cells.Item(4, 10).Value = 267297500
$ws.Cells.Item(4, 11).Value = 4847214.6
$ws.Cells.Item(4, 12).Value = 801892500
$ws.Cells.Item(4, 13).Value = -4847102.6
$ws.Cells.Item(4, 14).Value = -801892724

$ws.Cells.Item(6, 8).Value = 352.5
$ws.Cells.Item(6, 9).Value = 352.5
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 1057.5
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -944.5

$ws.Cells.Item(104, 8).Value = 7974.6665
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = 7974.6665
$ws.Cells.Item(104, 11).Value = 0
$ws.Cells.Item(104, 12).Value = 23923.9995
$ws.Cells.Item(104, 14).Value = -29165.9995

$ws.Cells.Item(109, 8).Value = 14142.091
$ws.Cells.Item(109, 9).Value = 1115.4
$ws.Cells.Item(109, 10).Value = 24997.666
$ws.Cells.Item(109, 11).Value = 3346.2
$ws.Cells.Item(109, 12).Value = 74992.99800000001
$ws.Cells.Item(109, 13).Value = -2306.2
$ws.Cells.Item(109, 14).Value = -77072.99800000001

$ws.Cells.Item(120, 8).Value = 21609.545
$ws.Cells.Item(120, 9).Value = 5541
$ws.Cells.Item(120, 10).Value = 35000
$ws.Cells.Item(120, 11).Value = 16623
$ws.Cells.Item(120, 12).Value = 105000
$ws.Cells.Item(120, 13).Value = -11785
$ws.Cells.Item(120, 14).Value = -114676

$ws.Cells.Item(122, 8).Value = 491.27274
$ws.Cells.Item(122, 9).Value = 496.66666
$ws.Cells.Item(122, 10).Value = 489.25
$ws.Cells.Item(122, 11).Value = 4469.99994
$ws.Cells.Item(122, 12).Value = 4403.25
$ws.Cells.Item(122, 13).Value = -2019.99994
$ws.Cells.Item(122, 14).Value = -9303.25

$ws.Cells.Item(138, 8).Value = 70542.336
$ws.Cells.Item(138, 9).Value = 146685.72
$ws.Cells.Item(138, 10).Value = 3916.875
$ws.Cells.Item(138, 11).Value = 440057.16
$ws.Cells.Item(138, 12).Value = 11750.625
$ws.Cells.Item(138, 13).Value = -434917.16
$ws.Cells.Item(138, 14).Value = -22030.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 3002.3333
$ws.Cells.Item(9, 9).Value = 2003.5
$ws.Cells.Item(9, 10).Value = 5000
$ws.Cells.Item(9, 11).Value = 2003.5
$ws.Cells.Item(9, 12).Value = 5000
$ws.Cells.Item(9, 13).Value = -1833.5
$ws.Cells.Item(9, 14).Value = -5340

$ws.Cells.Item(70, 8).Value = 9526908
$ws.Cells.Item(70, 9).Value = 15875183
$ws.Cells.Item(70, 10).Value = 4495
$ws.Cells.Item(70, 11).Value = 15875183
$ws.Cells.Item(70, 12).Value = 4495
$ws.Cells.Item(70, 13).Value = -15874913
$ws.Cells.Item(70, 14).Value = -5035

$ws.Cells.Item(73, 8).Value = 9526908
$ws.Cells.Item(73, 9).Value = 15875183
$ws.Cells.Item(73, 10).Value = 4495
$ws.Cells.Item(73, 11).Value = 15875183
$ws.Cells.Item(73, 12).Value = 4495
$ws.Cells.Item(73, 13).Value = -15874247
$ws.Cells.Item(73, 14).Value = -6367

$ws.Cells.Item(111, 8).Value = 30000
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = 30000
$ws.Cells.Item(111, 11).Value = 0
$ws.Cells.Item(111, 12).Value = 30000
$ws.Cells.Item(111, 14).Value = -36134

$ws.Cells.Item(122, 8).Value = 1836589.9
$ws.Cells.Item(122, 9).Value = 1836589.9
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 5509769.699999999
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -5507319.699999999

$ws.Cells.Item(126, 8).Value = 4726.0586
$ws.Cells.Item(126, 9).Value = 2171.4285
$ws.Cells.Item(126, 10).Value = 6514.3
$ws.Cells.Item(126, 11).Value = 6514.2855
$ws.Cells.Item(126, 12).Value = 19542.9
$ws.Cells.Item(126, 13).Value = -4044.2855
$ws.Cells.Item(126, 14).Value = -24482.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(114, 8).Value = 69999
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 69999
$ws.Cells.Item(114, 11).Value = 0
$ws.Cells.Item(114, 12).Value = 69999
$ws.Cells.Item(114, 14).Value = -78677

$ws.Cells.Item(122, 8).Value = 13277.667
$ws.Cells.Item(122, 9).Value = 5249.25
$ws.Cells.Item(122, 10).Value = 19700.4
$ws.Cells.Item(122, 11).Value = 15747.75
$ws.Cells.Item(122, 12).Value = 59101.2
$ws.Cells.Item(122, 13).Value = -13297.75
$ws.Cells.Item(122, 14).Value = -64001.2

$ws.Cells.Item(131, 8).Value = 74398.5
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 74398.5
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 74398.5
$ws.Cells.Item(131, 14).Value = -84478.5

$ws.Cells.Item(132, 8).Value = 3832.0833
$ws.Cells.Item(132, 9).Value = 2661.3333
$ws.Cells.Item(132, 10).Value = 6563.8335
$ws.Cells.Item(132, 11).Value = 7983.999899999999
$ws.Cells.Item(132, 12).Value = 19691.5005
$ws.Cells.Item(132, 13).Value = -5453.999899999999
$ws.Cells.Item(132, 14).Value = -24751.5005

$ws.Cells.Item(133, 8).Value = 67430
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 67430
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 67430
$ws.Cells.Item(133, 14).Value = -72490

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 3000
$ws.Cells.Item(3, 9).Value = 3000
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 3000
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = -2886

$ws.Cells.Item(70, 8).Value = 27198.334
$ws.Cells.Item(70, 9).Value = 27047.5
$ws.Cells.Item(70, 10).Value = 27500
$ws.Cells.Item(70, 11).Value = 27047.5
$ws.Cells.Item(70, 12).Value = 27500
$ws.Cells.Item(70, 13).Value = -26732.5
$ws.Cells.Item(70, 14).Value = -28130

$ws.Cells.Item(73, 8).Value = 27198.334
$ws.Cells.Item(73, 9).Value = 27047.5
$ws.Cells.Item(73, 10).Value = 27500
$ws.Cells.Item(73, 11).Value = 27047.5
$ws.Cells.Item(73, 12).Value = 27500
$ws.Cells.Item(73, 13).Value = -25955.5
$ws.Cells.Item(73, 14).Value = -29684

$ws.Cells.Item(122, 8).Value = 2487.3416
$ws.Cells.Item(122, 9).Value = 2233.743
$ws.Cells.Item(122, 10).Value = 3966.6667
$ws.Cells.Item(122, 11).Value = 6701.228999999999
$ws.Cells.Item(122, 12).Value = 11900.0001
$ws.Cells.Item(122, 13).Value = -4251.228999999999
$ws.Cells.Item(122, 14).Value = -16800.0001
